$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# ---------------------------------------------------------------------------
# Insert a new column "PopulationId" right after "IndividualId" (i.e. at C),
# shifting ModelParameterSheets..ModelFile one column to the right.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Header for the new column
$ws.Cells.Item(1, 3).Value2 = "PopulationId"
$ws.Cells.Item(1, 3).Font.Bold = $true

# ---------------------------------------------------------------------------
# Add a new row (row 4) describing a population-based scenario.
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value2 = "PopulationScenario"          # Scenario_name
$ws.Cells.Item(4, 2).Value2 = "Indiv"                        # IndividualId
$ws.Cells.Item(4, 3).Value2 = "TestPopulation"                # PopulationId
$ws.Cells.Item(4, 4).Value2 = "Global"                        # ModelParameterSheets
$ws.Cells.Item(4, 5).Value2 = "Aciclovir_iv_250mg"            # ApplicationProtocol
$ws.Cells.Item(4, 6).Value2 = 12                              # SimulationTime
$ws.Cells.Item(4, 7).Value2 = "h"                             # SimulationTimeUnit
$ws.Cells.Item(4, 8).Value2 = $false                          # SteadyState
# SteadyStateTime / SteadyStateTimeUnit (col 9/10) remain blank
$ws.Cells.Item(4, 11).Value2 = "Aciclovir.pkml"                # ModelFile

# ---------------------------------------------------------------------------
# Update the selection to match the authored state.
# ---------------------------------------------------------------------------
$ws.Range("I4:J4").Select() | Out-Null
